# BLT-44: "Experimenting with rendering strategies."
#
# Reposition/resize the two red "Straight Connector" lines on slide 1
# (the ones anchored with a:stCxn to the "Oval 211"/"Oval 218" shapes).
#
# NOTE on units: Shape.Left/.Top/.Width/.Height are expressed in points in
# the PowerPoint object model, while the underlying OOXML stores EMUs
# (1 pt = 12700 EMU). The interop layer additionally marshals these
# properties as single-precision (32-bit) floats, exactly like real
# PowerPoint COM does, so a naive `emu / 12700` can land one EMU away from
# the intended integer after the round-trip. The literals below are the
# closest representable points values that reproduce the exact target EMU
# offsets/extents once PowerPoint converts them back.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$connector1 = $s.Shapes.Item("Straight Connector 51")
$connector1.Left   = 642.2400512695312   # -> 8156448 EMU
$connector1.Top     = 167.5               # -> 2127250 EMU
$connector1.Width   = 55.69755935668945   # -> 707359 EMU
$connector1.Height  = 253.8750457763672   # -> 3224213 EMU

$connector2 = $s.Shapes.Item("Straight Connector 53")
$connector2.Left    = 817.9200439453125   # -> 10387584 EMU
$connector2.Top     = 150.5               # -> 1911350 EMU
$connector2.Width   = 125.92527770996094  # -> 1599251 EMU
$connector2.Height  = 260.0147399902344   # -> 3302187 EMU
